# 2007_to_2017_NAICS.xlsx - refresh crosswalk data + reposition selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the 2017 NAICS crosswalk code for 316999 (was 316990, should be 316998)
$ws.Range("B3").Value = 316998

# Leave the selection where the editor last left it before saving
$ws.Range("B4").Select()
